$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 36 (empty placeholder row for 44232 with no subscriber/price data).
#    Excel shifts all rows below it up by one and auto-adjusts formula references.
$ws.Rows(36).Delete()

# 2. Add the new "Sentiment" column header in J1.
$ws.Range("J1").Value = "Sentiment"

# 3. Fill J2:J61 with the new Sentiment data (aligned to the post-delete row numbering).
$sentiment = @(324,402,417,454,420,444,439,270,366,343,326,342,355,373,431,604,403,472,924,1411,1630,1159,1040,1211,2375,5882,3288,22488,77995,54281,38826,41411,18805,11433,2896,3314,3910,4336,2953,731,725,2211,1856,1604,1448,3056,7258,4267,1801,3073,2499,472,1804,3509,5672,2803,2320,1921,2345,1035)
$n = $sentiment.Length
$arr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $sentiment[$i]
}
$ws.Range("J2:J61").Value = $arr

# 4. Autofit the new column's width like Excel would after pasting in a new data column.
$ws.Columns("J:J").AutoFit()

# 5. Update the active selection to J2, matching the cell the user last touched.
$ws.Range("J2").Select()
